$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35, shifting existing rows 35:102 down to 36:103
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the new data record
$ws.Cells.Item(35, 1).Value = 10
$ws.Cells.Item(35, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(35, 3).Value = "La Araucanía"
$ws.Cells.Item(35, 4).Value = 44519
$ws.Cells.Item(35, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(35, 5).Value = 9
$ws.Cells.Item(35, 6).Value = 100112012
$ws.Cells.Item(35, 7).Value = "Espinaca"
$ws.Cells.Item(35, 8).Value = "Sin especificar"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 55
$ws.Cells.Item(35, 11).Value = 8000
$ws.Cells.Item(35, 12).Value = 8000
$ws.Cells.Item(35, 13).Value = 8000
$ws.Cells.Item(35, 14).Value = "`$/docena de atados"
$ws.Cells.Item(35, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(35, 16).Value = 2667
$ws.Cells.Item(35, 17).Value = 3
$ws.Cells.Item(35, 18).Value = "Hortaliza"
